$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Cntn2"
$ws.Cells.Item(2, 3).Value = "Cntnap2"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 1
$ws.Cells.Item(2, 6).Value = 0.3333333333333333
$ws.Cells.Item(2, 7).Value = 0.023556
$ws.Cells.Item(2, 8).Value = 0.070668
$ws.Cells.Item(2, 9).Value = 0.08088674619362546
$ws.Cells.Item(2, 10).Value = 0.08088674619362549
$ws.Cells.Item(2, 11).Value = 2
$ws.Cells.Item(2, 12).Value = 0.6666666666666666
$ws.Cells.Item(2, 13).Value = 0.001088
$ws.Cells.Item(2, 14).Value = 0.003264
$ws.Cells.Item(2, 15).Value = 0.02170357071613805
$ws.Cells.Item(2, 16).Value = 0.02170357071613804
$ws.Cells.Item(2, 17).Value = 0.000025628928
$ws.Cells.Item(2, 18).Value = 0.000230660352
$ws.Cells.Item(2, 19).Value = 0.00175553121601166
$ws.Cells.Item(2, 20).Value = 0.00175553121601166

# Row 3
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Cntn2"
$ws.Cells.Item(3, 3).Value = "Cntnap2"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 1
$ws.Cells.Item(3, 6).Value = 0.3333333333333333
$ws.Cells.Item(3, 7).Value = 0.023556
$ws.Cells.Item(3, 8).Value = 0.070668
$ws.Cells.Item(3, 9).Value = 0.08088674619362546
$ws.Cells.Item(3, 10).Value = 0.08088674619362549
$ws.Cells.Item(3, 11).Value = 2
$ws.Cells.Item(3, 12).Value = 0.6666666666666666
$ws.Cells.Item(3, 13).Value = 0.04018033333333333
$ws.Cells.Item(3, 14).Value = 0.120541
$ws.Cells.Item(3, 15).Value = 0.801522707626837
$ws.Cells.Item(3, 16).Value = 0.8015227076268369
$ws.Cells.Item(3, 17).Value = 0.0009464879319999998
$ws.Cells.Item(3, 18).Value = 0.008518391387999999
$ws.Cells.Item(3, 19).Value = 0.06483256382023943
$ws.Cells.Item(3, 20).Value = 0.06483256382023944

# Row 4
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Cntn2"
$ws.Cells.Item(4, 3).Value = "Cntnap2"
$ws.Cells.Item(4, 4).Value = "MuSCs"
$ws.Cells.Item(4, 5).Value = 1
$ws.Cells.Item(4, 6).Value = 0.3333333333333333
$ws.Cells.Item(4, 7).Value = 0.023556
$ws.Cells.Item(4, 8).Value = 0.070668
$ws.Cells.Item(4, 9).Value = 0.08088674619362546
$ws.Cells.Item(4, 10).Value = 0.08088674619362549
$ws.Cells.Item(4, 11).Value = 1
$ws.Cells.Item(4, 12).Value = 0.3333333333333333
$ws.Cells.Item(4, 13).Value = 0.008861666666666667
$ws.Cells.Item(4, 14).Value = 0.026585
$ws.Cells.Item(4, 15).Value = 0.1767737216570251
$ws.Cells.Item(4, 16).Value = 0.1767737216570251
$ws.Cells.Item(4, 17).Value = 0.00020874542
$ws.Cells.Item(4, 18).Value = 0.00187870878
$ws.Cells.Item(4, 19).Value = 0.01429865115737438
$ws.Cells.Item(4, 20).Value = 0.01429865115737438

# Row 5
$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "Cntn2"
$ws.Cells.Item(5, 3).Value = "Cntnap2"
$ws.Cells.Item(5, 4).Value = "ECs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 0.2501433333333333
$ws.Cells.Item(5, 8).Value = 0.7504299999999999
$ws.Cells.Item(5, 9).Value = 0.8589438069010353
$ws.Cells.Item(5, 10).Value = 0.8589438069010354
$ws.Cells.Item(5, 11).Value = 2
$ws.Cells.Item(5, 12).Value = 0.6666666666666666
$ws.Cells.Item(5, 13).Value = 0.001088
$ws.Cells.Item(5, 14).Value = 0.003264
$ws.Cells.Item(5, 15).Value = 0.02170357071613805
$ws.Cells.Item(5, 16).Value = 0.02170357071613804
$ws.Cells.Item(5, 17).Value = 0.0002721559466666667
$ws.Cells.Item(5, 18).Value = 0.00244940352
$ws.Cells.Item(5, 19).Value = 0.01864214765426544
$ws.Cells.Item(5, 20).Value = 0.01864214765426544

# Row 6
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Cntn2"
$ws.Cells.Item(6, 3).Value = "Cntnap2"
$ws.Cells.Item(6, 4).Value = "FAPs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 0.2501433333333333
$ws.Cells.Item(6, 8).Value = 0.7504299999999999
$ws.Cells.Item(6, 9).Value = 0.8589438069010353
$ws.Cells.Item(6, 10).Value = 0.8589438069010354
$ws.Cells.Item(6, 11).Value = 2
$ws.Cells.Item(6, 12).Value = 0.6666666666666666
$ws.Cells.Item(6, 13).Value = 0.04018033333333333
$ws.Cells.Item(6, 14).Value = 0.120541
$ws.Cells.Item(6, 15).Value = 0.801522707626837
$ws.Cells.Item(6, 16).Value = 0.8015227076268369
$ws.Cells.Item(6, 17).Value = 0.01005084251444444
$ws.Cells.Item(6, 18).Value = 0.09045758262999999
$ws.Cells.Item(6, 19).Value = 0.6884629658066208
$ws.Cells.Item(6, 20).Value = 0.6884629658066208

# Row 7
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Cntn2"
$ws.Cells.Item(7, 3).Value = "Cntnap2"
$ws.Cells.Item(7, 4).Value = "MuSCs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 0.2501433333333333
$ws.Cells.Item(7, 8).Value = 0.7504299999999999
$ws.Cells.Item(7, 9).Value = 0.8589438069010353
$ws.Cells.Item(7, 10).Value = 0.8589438069010354
$ws.Cells.Item(7, 11).Value = 1
$ws.Cells.Item(7, 12).Value = 0.3333333333333333
$ws.Cells.Item(7, 13).Value = 0.008861666666666667
$ws.Cells.Item(7, 14).Value = 0.026585
$ws.Cells.Item(7, 15).Value = 0.1767737216570251
$ws.Cells.Item(7, 16).Value = 0.1767737216570251
$ws.Cells.Item(7, 17).Value = 0.002216686838888889
$ws.Cells.Item(7, 18).Value = 0.01995018155
$ws.Cells.Item(7, 19).Value = 0.1518386934401491
$ws.Cells.Item(7, 20).Value = 0.1518386934401491

# Row 8
$ws.Cells.Item(8, 1).Value = "MuSCs"
$ws.Cells.Item(8, 2).Value = "Cntn2"
$ws.Cells.Item(8, 3).Value = "Cntnap2"
$ws.Cells.Item(8, 4).Value = "ECs"
$ws.Cells.Item(8, 5).Value = 1
$ws.Cells.Item(8, 6).Value = 0.3333333333333333
$ws.Cells.Item(8, 7).Value = 0.01752266666666667
$ws.Cells.Item(8, 8).Value = 0.052568
$ws.Cells.Item(8, 9).Value = 0.0601694469053391
$ws.Cells.Item(8, 10).Value = 0.06016944690533912
$ws.Cells.Item(8, 11).Value = 2
$ws.Cells.Item(8, 12).Value = 0.6666666666666666
$ws.Cells.Item(8, 13).Value = 0.001088
$ws.Cells.Item(8, 14).Value = 0.003264
$ws.Cells.Item(8, 15).Value = 0.02170357071613805
$ws.Cells.Item(8, 16).Value = 0.02170357071613804
$ws.Cells.Item(8, 17).Value = 0.00001906466133333334
$ws.Cells.Item(8, 18).Value = 0.000171581952
$ws.Cells.Item(8, 19).Value = 0.001305891845860941
$ws.Cells.Item(8, 20).Value = 0.001305891845860941

# Row 9
$ws.Cells.Item(9, 1).Value = "MuSCs"
$ws.Cells.Item(9, 2).Value = "Cntn2"
$ws.Cells.Item(9, 3).Value = "Cntnap2"
$ws.Cells.Item(9, 4).Value = "FAPs"
$ws.Cells.Item(9, 5).Value = 1
$ws.Cells.Item(9, 6).Value = 0.3333333333333333
$ws.Cells.Item(9, 7).Value = 0.01752266666666667
$ws.Cells.Item(9, 8).Value = 0.052568
$ws.Cells.Item(9, 9).Value = 0.0601694469053391
$ws.Cells.Item(9, 10).Value = 0.06016944690533912
$ws.Cells.Item(9, 11).Value = 2
$ws.Cells.Item(9, 12).Value = 0.6666666666666666
$ws.Cells.Item(9, 13).Value = 0.04018033333333333
$ws.Cells.Item(9, 14).Value = 0.120541
$ws.Cells.Item(9, 15).Value = 0.801522707626837
$ws.Cells.Item(9, 16).Value = 0.8015227076268369
$ws.Cells.Item(9, 17).Value = 0.0007040665875555555
$ws.Cells.Item(9, 18).Value = 0.006336599288
$ws.Cells.Item(9, 19).Value = 0.0482271779999766
$ws.Cells.Item(9, 20).Value = 0.04822717799997661

# Row 10
$ws.Cells.Item(10, 1).Value = "MuSCs"
$ws.Cells.Item(10, 2).Value = "Cntn2"
$ws.Cells.Item(10, 3).Value = "Cntnap2"
$ws.Cells.Item(10, 4).Value = "MuSCs"
$ws.Cells.Item(10, 5).Value = 1
$ws.Cells.Item(10, 6).Value = 0.3333333333333333
$ws.Cells.Item(10, 7).Value = 0.01752266666666667
$ws.Cells.Item(10, 8).Value = 0.052568
$ws.Cells.Item(10, 9).Value = 0.0601694469053391
$ws.Cells.Item(10, 10).Value = 0.06016944690533912
$ws.Cells.Item(10, 11).Value = 1
$ws.Cells.Item(10, 12).Value = 0.3333333333333333
$ws.Cells.Item(10, 13).Value = 0.008861666666666667
$ws.Cells.Item(10, 14).Value = 0.026585
$ws.Cells.Item(10, 15).Value = 0.1767737216570251
$ws.Cells.Item(10, 16).Value = 0.1767737216570251
$ws.Cells.Item(10, 17).Value = 0.0001552800311111111
$ws.Cells.Item(10, 18).Value = 0.00139752028
$ws.Cells.Item(10, 19).Value = 0.01063637705950156
$ws.Cells.Item(10, 20).Value = 0.01063637705950157
